$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15 (pushes existing rows 15-66 down to 16-67)
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly data point
$ws.Cells.Item(15, 1).Value  = 7
$ws.Cells.Item(15, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(15, 3).Value  = "Ñuble"
$ws.Cells.Item(15, 4).Value  = 44565
$ws.Cells.Item(15, 5).Value  = 16
$ws.Cells.Item(15, 6).Value  = 100112030
$ws.Cells.Item(15, 7).Value  = "Poroto granado"
$ws.Cells.Item(15, 8).Value  = "Sin especificar"
$ws.Cells.Item(15, 9).Value  = "Primera"
$ws.Cells.Item(15, 10).Value = 120
$ws.Cells.Item(15, 11).Value = 31000
$ws.Cells.Item(15, 12).Value = 32000
$ws.Cells.Item(15, 13).Value = 31500
$ws.Cells.Item(15, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Región del Maule"
$ws.Cells.Item(15, 16).Value = 1260
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# Match the date cell formatting used by the rest of column D
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
